# cypress/fixtures/customers.xlsx
# feat: Enable AgGrid column grouping, and add import/export in customToolbar
#
# - Update Howard's birthday (customers!E5)
# - Add a new customer row (Billy) to the customers sheet
# - Make the "customers" sheet the active/selected sheet (was "UnusedSheet")

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # customers

# Howard's birthday text changed
$ws1.Range("E5").Value = "21/05/2002"

# New row for Billy
$ws1.Range("A6").Value = "Billy"
$ws1.Range("B6").Value = 23

$ws1.Range("C6").Value = $true
$ws1.Range("C6").NumberFormat = """TRUE"";""TRUE"";""FALSE"""

$ws1.Range("D6").Value = "Beer"

$ws1.Range("E6").NumberFormat = "mm/dd/yy"
$ws1.Range("E6").Formula = "=DATE(1940,4,28)"

$ws1.Range("F6").Value = 1.25

# customers becomes the active sheet / tab, with E7 selected
$ws1.Activate()
$ws1.Range("E7").Select() | Out-Null
